$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Add a new translation row ("message queue for temperature") above the
# existing SingleUseId2 row, shifting it down from row 5 to row 6.
$ws.Rows("5:5").Insert()

$ws.Range("B5").Value = "ResourceId1"
$ws.Range("C5").Value = "Large"
$ws.Range("D5").Value = "Left"
$ws.Range("E5").Value = "LTR"

# F5 needs to hold the literal text "0123456789." (a numeric wildcard mask
# with a leading zero and trailing dot). Assigning that string straight to
# .Value would let Excel's "looks like a number" auto-detection convert it
# to 123456789 and drop the formatting-significant characters. Instead,
# stage the exact text as a formula result in a scratch cell far outside
# the sheet's used range, copy its computed value (which keeps the Text
# type without touching the cell's style), paste that into F5, then clean
# the scratch cell back up.
$ws.Range("Z100").Formula = "=""0123456789."""
$excel.Calculate()
$ws.Range("Z100").Copy()
$ws.Range("F5").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
$ws.Range("Z100").ClearContents()
